# Fix the "Skor" (Score) column so every row shows the full scoring
# criteria text, and wrap that text so it is fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skorText = "Skor`nAda dan berfungsi semuanya : 2 Ada tapi tidak berfungsi sebagian : 1`nAda dan tidak berfungsi semuanya : 0"

# E4 previously just said "Skor" while the rest of the column (E2,E3,E5,E6,E7)
# had the full multi-line scoring criteria text - make it consistent.
$ws.Range("E4").Value = $skorText

# Wrap text in the whole "Skor" column so the multi-line criteria is readable.
$ws.Range("E2:E7").WrapText = $true

# Keep the original row heights (wrapping would otherwise auto-grow row 4).
$ws.Rows("2").RowHeight = 17.25
$ws.Rows("3").RowHeight = 17.25
$ws.Rows("4").RowHeight = 17.25
$ws.Rows("5").RowHeight = 17.25
$ws.Rows("6").RowHeight = 17.25
$ws.Rows("7").RowHeight = 18
